$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: "Material Type" - bold white text on a black fill, centered.
$header = $ws.Range("G1")
$header.Value = "Material Type"
$header.Font.Bold = $true
$header.Font.Color = 16777215
$header.Interior.Color = 0
$header.Interior.PatternColor = 0
$header.HorizontalAlignment = -4108

# Data cells: "DNA:Genomic" for every data row, centered (same look as other data columns).
$dataRange = $ws.Range("G2:G24")
$dataRange.Value = "DNA:Genomic"
$dataRange.HorizontalAlignment = -4108

# Keep the sheet selection consistent with the newly populated column.
$ws.Range("G1:G24").Select()
